# Update the evaluation figures (Note_EPP, MNG, Facteur) for the two
# teams on the "Sommaire de l'EPP" sheet.
#
# Column E = Note_EPP (raw score)
# Column F = MNG      (team average of Note_EPP)
# Column G = Facteur  (Note_EPP / MNG)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team 1 (rows 2-5) ---------------------------------------------------
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 46.25
$ws.Range("G2").Value = 1.297297297297297

$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 46.25
$ws.Range("G3").Value = 1.621621621621622

$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 46.25
$ws.Range("G4").Value = 0.5405405405405406

$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 46.25
$ws.Range("G5").Value = 0.5405405405405406

# --- Team 2 (rows 6-9) ----------------------------------------------------
$ws.Range("E6").Value = 66.66666666666667
$ws.Range("F6").Value = 50.00000000000001
$ws.Range("G6").Value = 1.333333333333333

$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 50.00000000000001
$ws.Range("G7").Value = 2

$ws.Range("E8").Value = 33.33333333333334
$ws.Range("F8").Value = 50.00000000000001
$ws.Range("G8").Value = 0.6666666666666666

$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 50.00000000000001
$ws.Range("G9").Value = 0
